$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("db_table")
$ws.Activate()

# Clear the contents of row 18 (A18:BE18), keeping the existing cell
# formatting/styles in place - this mirrors selecting the row and
# pressing Delete in the Excel UI.
$ws.Range("A18:BE18").ClearContents()

# Update the frozen-pane anchor and the active selection to reflect
# where the user left off after clearing the row.
$ws.Range("S3").Select()
$ws.Application.ActiveWindow.FreezePanes = $true

$ws.Rows.Item(18).Select()
